# This script updates the "cryptos" price-tracker worksheet in place,
# applying the latest price/volume snapshot (refreshed by the scheduled
# GitHub Actions job) to rows 2-51 of the active sheet.
#
# Price values in column D are plain text (the source data uses "."
# as a thousands separator in several rows, e.g. "27.605.85"), so any
# value that Excel would otherwise auto-convert to a number is written
# with a leading apostrophe to force a text cell, and the style is reset
# back to "Normal" immediately after so no stray formatting is left
# behind (Excel applies a "Text" number format to quote-prefixed cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.605.85"
$ws.Range("E2").Value = "  +0.41%  "
# Row 3
$ws.Range("D3").Value = "1.845.87"
$ws.Range("E3").Value = "  +0.23%  "
# Row 4
$ws.Range("D4").Value = "'1.029"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
$ws.Range("D5").Value = "'321.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
# Row 6
$ws.Range("D6").Value = "'1.026"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.12%  "
# Row 7
$ws.Range("D7").Value = "'0.4372"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.07%  "
# Row 8
$ws.Range("D8").Value = "'0.3782"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.38%  "
# Row 9
$ws.Range("D9").Value = "'0.07368"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.25%  "
# Row 10
$ws.Range("D10").Value = "'0.8804"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "
# Row 11
$ws.Range("D11").Value = "'21.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
# Row 12
$ws.Range("D12").Value = "1.852.04"
$ws.Range("E12").Value = "  -0.85%  "
# Row 13
$ws.Range("D13").Value = "'5.486"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
# Row 14
$ws.Range("D14").Value = "'6.702"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
# Row 15
$ws.Range("D15").Value = "'0.07127"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.27%  "
# Row 16
$ws.Range("D16").Value = "'85.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.93%  "
# Row 17
$ws.Range("D17").Value = "'1.031"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
# Row 18
$ws.Range("D18").Value = "'0.000009044"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.17%  "
# Row 19
$ws.Range("E19").Value = "  -0.09%  "
# Row 20
$ws.Range("D20").Value = "'15.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
# Row 21
$ws.Range("D21").Value = "27.606.96"
$ws.Range("E21").Value = "  +0.33%  "
# Row 22
$ws.Range("D22").Value = "'5.282"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.07%  "
# Row 23
$ws.Range("D23").Value = "'11.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.65%  "
# Row 24
$ws.Range("D24").Value = "2.084.96"
$ws.Range("E24").Value = "  +0.29%  "
# Row 25
$ws.Range("D25").Value = "'2.026"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.36%  "
# Row 26
$ws.Range("D26").Value = "'157.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
# Row 27
$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.16%  "
# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'1.987"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.60%  "
# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'5.322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.40%  "
# Row 30
$ws.Range("D30").Value = "'117.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.25%  "
# Row 31
$ws.Range("D31").Value = "'0.08996"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "
# Row 32
$ws.Range("D32").Value = "'0.7710"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.20%  "
# Row 33
$ws.Range("D33").Value = "'1.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.17%  "
# Row 34
$ws.Range("D34").Value = "'2.989"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.96%  "
# Row 35
$ws.Range("D35").Value = "'4.543"
$ws.Range("D35").Style = "Normal"
# Row 36
$ws.Range("D36").Value = "'1.027"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.21%  "
# Row 37
$ws.Range("E37").Value = "  -0.65%  "
# Row 38
$ws.Range("D38").Value = "'0.01968"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.09%  "
# Row 39
$ws.Range("D39").Value = "'0.05263"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.20%  "
# Row 40
$ws.Range("D40").Value = "'2.839"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.08%  "
# Row 41
$ws.Range("D41").Value = "'0.5166"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
# Row 42
$ws.Range("D42").Value = "'0.1668"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "
# Row 43
$ws.Range("D43").Value = "'6.832"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.09%  "
# Row 44
$ws.Range("D44").Value = "'8.774"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "
# Row 45
$ws.Range("D45").Value = "'109.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.71%  "
# Row 46
$ws.Range("D46").Value = "'10.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.65%  "
# Row 47
$ws.Range("D47").Value = "'0.06592"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.17%  "
# Row 48
$ws.Range("D48").Value = "'1.028"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
# Row 49
$ws.Range("D49").Value = "'1.697"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.49%  "
# Row 50
$ws.Range("D50").Value = "'0.4683"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "
# Row 51
$ws.Range("D51").Value = "'1.890"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
